$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell C2 value from 6.4 to 5
$ws.Range("C2").Value = 5

# Update the selection from C6 to C4
$ws.Range("C4").Select()
